$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B3 from 1 to 0.5, and C3 from 2 to 1.
# C2 contains formula =B3 and will recalc to 0.5 automatically.
$ws.Range("B3").Value = 0.5
$ws.Range("C3").Value = 1

$wb.Application.Calculate()
